$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNOW")

$ws.Range("B12").Value = 6000000.0
$ws.Range("C12").Value = 5000000.0
$ws.Range("D12").Value = 8000000.0
$ws.Range("E12").Value = 6000000.0
